$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure these cells keep being stored as plain text (inline strings),
# matching how the dates are currently stored in the sheet, rather than
# being auto-converted to Excel date serial numbers.
$cells = @("D6","D7","D8","D9","C17","D18","D19","D21","D22")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D6").Value = "2024-06-30"
$ws.Range("D7").Value = "2024-06-30"
$ws.Range("D8").Value = "2025-01-19"
$ws.Range("D9").Value = "2024-06-30"
$ws.Range("C17").Value = "1992-09-26"
$ws.Range("D18").Value = "2024-12-31"
$ws.Range("D19").Value = "2023-06-07"
$ws.Range("D21").Value = "2025-01-01"
$ws.Range("D22").Value = "2023-06-07"
